# Recompute ligand-receptor (TPM) derived metrics for Lama2-Rpsa sheet
# after updated TPM normalization values (commit: "update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.667069666666666
$ws.Range("H2").Value = 11.001209
$ws.Range("I2").Value = 0.01298011522000835
$ws.Range("J2").Value = 0.01298011522000835
$ws.Range("M2").Value = 68.18146900000001
$ws.Range("N2").Value = 204.544407
$ws.Range("O2").Value = 0.1244286043321187
$ws.Range("P2").Value = 0.1244286043321187
$ws.Range("Q2").Value = 250.0261967986737
$ws.Range("R2").Value = 2250.235771188063
$ws.Range("S2").Value = 0.001615097620895731
$ws.Range("T2").Value = 0.001615097620895731
$ws.Range("G3").Value = 3.667069666666666
$ws.Range("H3").Value = 11.001209
$ws.Range("I3").Value = 0.01298011522000835
$ws.Range("J3").Value = 0.01298011522000835
$ws.Range("O3").Value = 0.345973452289334
$ws.Range("P3").Value = 0.3459734522893341
$ws.Range("Q3").Value = 695.1972734365932
$ws.Range("R3").Value = 6256.775460929339
$ws.Range("S3").Value = 0.004490775273779618
$ws.Range("T3").Value = 0.004490775273779619
$ws.Range("G4").Value = 3.667069666666666
$ws.Range("H4").Value = 11.001209
$ws.Range("I4").Value = 0.01298011522000835
$ws.Range("J4").Value = 0.01298011522000835
$ws.Range("M4").Value = 188.0130056666667
$ws.Range("N4").Value = 564.0390170000001
$ws.Range("O4").Value = 0.3431166302883566
$ws.Range("P4").Value = 0.3431166302883567
$ws.Range("Q4").Value = 689.4567900190615
$ws.Range("R4").Value = 6205.111110171553
$ws.Range("S4").Value = 0.004453693395043877
$ws.Range("T4").Value = 0.004453693395043879
$ws.Range("G5").Value = 3.667069666666666
$ws.Range("H5").Value = 11.001209
$ws.Range("I5").Value = 0.01298011522000835
$ws.Range("J5").Value = 0.01298011522000835
$ws.Range("M5").Value = 102.1836573333333
$ws.Range("N5").Value = 306.550972
$ws.Range("O5").Value = 0.1864813130901906
$ws.Range("P5").Value = 0.1864813130901907
$ws.Range("Q5").Value = 374.7145902361275
$ws.Range("R5").Value = 3372.431312125148
$ws.Range("S5").Value = 0.002420548930289126
$ws.Range("T5").Value = 0.002420548930289127
$ws.Range("I6").Value = 0.5954329572989919
$ws.Range("J6").Value = 0.595432957298992
$ws.Range("M6").Value = 68.18146900000001
$ws.Range("N6").Value = 204.544407
$ws.Range("O6").Value = 0.1244286043321187
$ws.Range("P6").Value = 0.1244286043321187
$ws.Range("Q6").Value = 11469.37721574079
$ws.Range("R6").Value = 103224.3949416671
$ws.Range("S6").Value = 0.07408889185005957
$ws.Range("T6").Value = 0.07408889185005958
$ws.Range("I7").Value = 0.5954329572989919
$ws.Range("J7").Value = 0.595432957298992
$ws.Range("O7").Value = 0.345973452289334
$ws.Range("P7").Value = 0.3459734522893341
$ws.Range("S7").Value = 0.2060039958435798
$ws.Range("T7").Value = 0.2060039958435799
$ws.Range("I8").Value = 0.5954329572989919
$ws.Range("J8").Value = 0.595432957298992
$ws.Range("M8").Value = 188.0130056666667
$ws.Range("N8").Value = 564.0390170000001
$ws.Range("O8").Value = 0.3431166302883566
$ws.Range("P8").Value = 0.3431166302883567
$ws.Range("Q8").Value = 31627.24586435957
$ws.Range("R8").Value = 284645.2127792361
$ws.Range("S8").Value = 0.2043029498710611
$ws.Range("T8").Value = 0.2043029498710611
$ws.Range("I9").Value = 0.5954329572989919
$ws.Range("J9").Value = 0.595432957298992
$ws.Range("M9").Value = 102.1836573333333
$ws.Range("N9").Value = 306.550972
$ws.Range("O9").Value = 0.1864813130901906
$ws.Range("P9").Value = 0.1864813130901907
$ws.Range("Q9").Value = 17189.17072965966
$ws.Range("R9").Value = 154702.5365669369
$ws.Range("S9").Value = 0.1110371197342914
$ws.Range("T9").Value = 0.1110371197342915
$ws.Range("G10").Value = 110.4727123333333
$ws.Range("H10").Value = 331.418137
$ws.Range("I10").Value = 0.3910338949346852
$ws.Range("J10").Value = 0.3910338949346853
$ws.Range("M10").Value = 68.18146900000001
$ws.Range("N10").Value = 204.544407
$ws.Range("O10").Value = 0.1244286043321187
$ws.Range("P10").Value = 0.1244286043321187
$ws.Range("Q10").Value = 7532.191811301085
$ws.Range("R10").Value = 67789.72630170976
$ws.Range("S10").Value = 0.04865580179327521
$ws.Range("T10").Value = 0.04865580179327521
$ws.Range("G11").Value = 110.4727123333333
$ws.Range("H11").Value = 331.418137
$ws.Range("I11").Value = 0.3910338949346852
$ws.Range("J11").Value = 0.3910338949346853
$ws.Range("O11").Value = 0.345973452289334
$ws.Range("P11").Value = 0.3459734522893341
$ws.Range("Q11").Value = 20943.24225726784
$ws.Range("R11").Value = 188489.1803154106
$ws.Range("S11").Value = 0.1352873465926978
$ws.Range("T11").Value = 0.1352873465926978
$ws.Range("G12").Value = 110.4727123333333
$ws.Range("H12").Value = 331.418137
$ws.Range("I12").Value = 0.3910338949346852
$ws.Range("J12").Value = 0.3910338949346853
$ws.Range("M12").Value = 188.0130056666667
$ws.Range("N12").Value = 564.0390170000001
$ws.Range("O12").Value = 0.3431166302883566
$ws.Range("P12").Value = 0.3431166302883567
$ws.Range("Q12").Value = 20770.30668993904
$ws.Range("R12").Value = 186932.7602094513
$ws.Range("S12").Value = 0.1341702323585205
$ws.Range("T12").Value = 0.1341702323585205
$ws.Range("G13").Value = 110.4727123333333
$ws.Range("H13").Value = 331.418137
$ws.Range("I13").Value = 0.3910338949346852
$ws.Range("J13").Value = 0.3910338949346853
$ws.Range("M13").Value = 102.1836573333333
$ws.Range("N13").Value = 306.550972
$ws.Range("O13").Value = 0.1864813130901906
$ws.Range("P13").Value = 0.1864813130901907
$ws.Range("Q13").Value = 11288.50578175324
$ws.Range("R13").Value = 101596.5520357792
$ws.Range("S13").Value = 0.07292051419019174
$ws.Range("T13").Value = 0.07292051419019177
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.1562396666666667
$ws.Range("H14").Value = 0.468719
$ws.Range("I14").Value = 0.0005530325463144183
$ws.Range("J14").Value = 0.0005530325463144184
$ws.Range("M14").Value = 68.18146900000001
$ws.Range("N14").Value = 204.544407
$ws.Range("O14").Value = 0.1244286043321187
$ws.Range("P14").Value = 0.1244286043321187
$ws.Range("Q14").Value = 10.65264998940367
$ws.Range("R14").Value = 95.873849904633
$ws.Range("S14").Value = 0.00006881306788814084
$ws.Range("T14").Value = 0.00006881306788814085
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.1562396666666667
$ws.Range("H15").Value = 0.468719
$ws.Range("I15").Value = 0.0005530325463144183
$ws.Range("J15").Value = 0.0005530325463144184
$ws.Range("O15").Value = 0.345973452289334
$ws.Range("P15").Value = 0.3459734522893341
$ws.Range("Q15").Value = 29.61966914799333
$ws.Range("R15").Value = 266.57702233194
$ws.Range("S15").Value = 0.0001913345792767603
$ws.Range("T15").Value = 0.0001913345792767604
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.1562396666666667
$ws.Range("H16").Value = 0.468719
$ws.Range("I16").Value = 0.0005530325463144183
$ws.Range("J16").Value = 0.0005530325463144184
$ws.Range("M16").Value = 188.0130056666667
$ws.Range("N16").Value = 564.0390170000001
$ws.Range("O16").Value = 0.3431166302883566
$ws.Range("P16").Value = 0.3431166302883567
$ws.Range("Q16").Value = 29.37508933435812
$ws.Range("R16").Value = 264.375804009223
$ws.Range("S16").Value = 0.0001897546637311927
$ws.Range("T16").Value = 0.0001897546637311928
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.1562396666666667
$ws.Range("H17").Value = 0.468719
$ws.Range("I17").Value = 0.0005530325463144183
$ws.Range("J17").Value = 0.0005530325463144184
$ws.Range("M17").Value = 102.1836573333333
$ws.Range("N17").Value = 306.550972
$ws.Range("O17").Value = 0.1864813130901906
$ws.Range("P17").Value = 0.1864813130901907
$ws.Range("Q17").Value = 15.96514056054089
$ws.Range("R17").Value = 143.686265044868
$ws.Range("S17").Value = 0.0001031302354183244
$ws.Range("T17").Value = 0.0001031302354183244
